$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top to hold header labels, shifting all data down by one row.
$ws.Rows.Item(1).Insert()

# Set header values in the newly inserted row 1
$ws.Cells.Item(1, 1).Value = "dato"
$ws.Cells.Item(1, 2).Value = "styringsrent PPR 1/23"

# Set column B width as in the target file (OOXML width=18.5 -> ColumnWidth offset by 5/6)
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668

# Update selection to match target (active cell C8)
$ws.Range("C8").Select()
